$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cover sheet")
Write-Host "Sheet1 name: $($ws.Name)"
